# "Generate Report for Handback" -------------------------------------------
# Localization CI report: the handback pass completed, so:
#   * the report Status flips from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   * each locale sheet's "Latest Target File" / "Latest Handback File" /
#     "Latest Handback DateTime" columns get filled in with the handback
#     artifact info (file names, timestamps, and a hyperlink to the source)
#   * a few columns are widened so the new long file names / links fit

$wb  = $excel.ActiveWorkbook
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

$srcMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/338599ba1638b116fcc0e5b4fec7af0337dda4b3/e2e/2bf6a623-283a-48ef-9d6e-1f7852ead1d9.md"
$srcMdDisplay = "2bf6a623-283a-48ef-9d6e-1f7852ead1d9.md"

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US"-
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn: Latest Target File (I) + Latest Handback File (J) -------------
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $srcMdUrl, "", "", $srcMdDisplay)
$zhcn.Range("J2").Value = "2bf6a623-283a-48ef-9d6e-1f7852ead1d9.620d4ce6edac63726341fd8012fba993f62faff2.zh-cn.xlf"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $srcMdUrl, "", "", $srcMdDisplay)
$zhcn.Range("J3").Value = "2bf6a623-283a-48ef-9d6e-1f7852ead1d9.620d4ce6edac63726341fd8012fba993f62faff2.zh-cn.xlf"

# --- de-de: Latest Target File (I) + Latest Handback File (J) + ----------
#     Latest Handback DateTime (K)
$dede.Hyperlinks.Add($dede.Range("I2"), $srcMdUrl, "", "", $srcMdDisplay)
$dede.Range("J2").Value = "2bf6a623-283a-48ef-9d6e-1f7852ead1d9.620d4ce6edac63726341fd8012fba993f62faff2.de-de.xlf"
$dede.Range("K2").Value = "2016-08-24 19:12:54"

$dede.Hyperlinks.Add($dede.Range("I3"), $srcMdUrl, "", "", $srcMdDisplay)
$dede.Range("J3").Value = "2bf6a623-283a-48ef-9d6e-1f7852ead1d9.620d4ce6edac63726341fd8012fba993f62faff2.de-de.xlf"
$dede.Range("K3").Value = "2016-08-24 19:12:54"

# --- Column width adjustments ---------------------------------------------
# Overview: zh-cn / de-de summary columns (E, F) widen
$overview.Range("E1").ColumnWidth = 29.166666666666668
$overview.Range("F1").ColumnWidth = 29.166666666666668

# zh-cn / de-de: Status (C) widens; Latest Target File (I) / Latest
# Handback File (J) widen to fit the new link text / file names
$zhcn.Range("C1").ColumnWidth = 29.166666666666668
$zhcn.Range("I1").ColumnWidth = 39.166666666666664
$zhcn.Range("J1").ColumnWidth = 39.166666666666664

$dede.Range("C1").ColumnWidth = 29.166666666666668
$dede.Range("I1").ColumnWidth = 39.166666666666664
$dede.Range("J1").ColumnWidth = 39.166666666666664
